# Fills in the four "half-year summary" comment/grade cells that were left
# empty (<w:t xml:space="preserve"/>) in the report-card tables for Naomi Jung.
#
# Each of the 4 subject tables (English, Arithmetic, Torah, Prophets) has the
# same 2x2 (merged) layout:
#   row1 col1 = subject name   | row1 col2 = free-text comment  (vMerge restart)
#   row2 col1 = "ciun:" label + a 2nd paragraph holding the grade number
#
# $d.Paragraphs.Item(n) walks every paragraph in the document, including the
# ones nested inside table cells, in document order. Cell-relative paragraph
# indexing (Cell.Range.Paragraphs.Item(2)) collapses onto paragraph 1 for this
# 2-paragraph cell in this host, so we address the target paragraphs by their
# absolute position in $d.Paragraphs instead (verified against the raw XML).

$d = $word.ActiveDocument
$nl = [char]10   # literal line-feed kept *inside* a single <w:t> run (not a new <w:p>)

# Table 1 (אנגלית / English)   - comment
$d.Paragraphs.Item(5).Range.Text = "במחצית זאת התמקדנו על האותיות,שיננו את ההברות שכל אות משמיע,והתחלנו לעבוד על קריאה בסיסית," + $nl + "נעמי את ילדה משקיעה ואכפתית, בהצלחה רבה!"

# Table 1 (אנגלית / English)   - grade
$d.Paragraphs.Item(8).Range.Text = "92"

# Table 2 (חשבון / Arithmetic) - comment
$d.Paragraphs.Item(15).Range.Text = "במחצית זאת למדנו כפל במאונך,חילוק ארוך, הרחבנו את סדר הפעולות והשימוש בסוגריים, וכן את התכונות ב0 ו1." + $nl + "וכן למדנו את עיקרון השבר הפשוט." + $nl + "נעמי את ילדה מעולה, עלי והצליחי!"

# Table 2 (חשבון / Arithmetic) - grade
$d.Paragraphs.Item(18).Range.Text = "87"

# Table 3 (תורה / Torah)       - comment
$d.Paragraphs.Item(25).Range.Text = "במחצית זאת למדנו את חומש `"במדבר`", למדנו על מסעות עם ישראל והפקנו לקחים רבים מכך!" + $nl + "נעמי את תלמידה מעולה, הרבה בהצלחה!"

# Table 3 (תורה / Torah)       - grade
$d.Paragraphs.Item(28).Range.Text = "98"

# Table 4 (נביא / Prophets)    - comment
$d.Paragraphs.Item(35).Range.Text = "במחצית זאת למדנו את ספר נביא `"יהושע`", למדנו ועקבנו אחרי פעולתיו ומעשיו בחייו ובדורו, וקיבלנו מסרים רבים!" + $nl + "נעמי את תלמידה מעולה, עלי והצלחי!"

# Table 4 (נביא / Prophets)    - grade
$d.Paragraphs.Item(38).Range.Text = "84"
